$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.198912084953569
$ws.Range("D2").Value = 8.974613108529622
$ws.Range("E2").Value = 12.9751858879603
$ws.Range("F2").Value = 32.46197151789978
$ws.Range("G2").Value = 3.628922753995076
$ws.Range("I2").Value = 17.91155213406074
$ws.Range("J2").Value = 9.647527558444178
$ws.Range("M2").Value = 20.19054413580184
$ws.Range("N2").Value = 17.23930263989935
$ws.Range("O2").Value = 24.55023775341948
$ws.Range("B3").Value = 6.997507437259769
$ws.Range("D3").Value = 8.998735157011918
$ws.Range("E3").Value = 13.03237523884615
$ws.Range("F3").Value = 32.29428320408314
$ws.Range("G3").Value = 3.631934186810442
$ws.Range("I3").Value = 17.95571653231028
$ws.Range("J3").Value = 9.684353539510559
$ws.Range("M3").Value = 19.72260192855842
$ws.Range("N3").Value = 17.18880106548365
$ws.Range("O3").Value = 24.42130678895314
$ws.Range("B4").Value = 6.871288594024167
$ws.Range("D4").Value = 9.014477137057479
$ws.Range("E4").Value = 13.06944165070674
$ws.Range("F4").Value = 32.20063308323169
$ws.Range("G4").Value = 3.63388164954294
$ws.Range("I4").Value = 17.98659562100642
$ws.Range("J4").Value = 9.708137620153819
$ws.Range("M4").Value = 19.4313280720625
$ws.Range("N4").Value = 17.15975828268506
$ws.Range("O4").Value = 24.34890014845619
$ws.Range("B5").Value = 6.819282695550045
$ws.Range("D5").Value = 9.021126592267608
$ws.Range("E5").Value = 13.08503808661739
$ws.Range("F5").Value = 32.16484081968659
$ws.Range("G5").Value = 3.634700094328831
$ws.Range("I5").Value = 18.00012564123216
$ws.Range("J5").Value = 9.718125425726909
$ws.Range("M5").Value = 19.31179374640394
$ws.Range("N5").Value = 17.14842603599106
$ws.Range("O5").Value = 24.32111591738361
$ws.Range("B6").Value = 6.810614883533242
$ws.Range("D6").Value = 9.022244901847579
$ws.Range("E6").Value = 13.08765757489094
$ws.Range("F6").Value = 32.159041525458
$ws.Range("G6").Value = 3.634837499159083
$ws.Range("I6").Value = 18.00242947341246
$ws.Range("J6").Value = 9.719801767226933
$ws.Range("M6").Value = 19.29189927499916
$ws.Range("N6").Value = 17.14657494776142
$ws.Range("O6").Value = 24.31660698691072
$ws.Range("B7").Value = 6.870589436429229
$ws.Range("D7").Value = 9.014565864062407
$ws.Range("E7").Value = 13.06964999834189
$ws.Range("F7").Value = 32.20014074048311
$ws.Range("G7").Value = 3.633892586692688
$ws.Range("I7").Value = 17.986774258427
$ws.Range("J7").Value = 9.708271121403707
$ws.Range("M7").Value = 19.42971917063073
$ws.Range("N7").Value = 17.15960340449234
$ws.Range("O7").Value = 24.34851844006556
$ws.Range("B8").Value = 7.130039814457801
$ws.Range("D8").Value = 8.982737399925057
$ws.Range("E8").Value = 12.99450010386875
$ws.Range("F8").Value = 32.40223896520807
$ws.Range("G8").Value = 3.629940718969742
$ws.Range("I8").Value = 17.92599959778081
$ws.Range("J8").Value = 9.659982165441161
$ws.Range("M8").Value = 20.03011044795811
$ws.Range("N8").Value = 17.22148710960409
$ws.Range("O8").Value = 24.50439439249589
$ws.Range("B9").Value = 7.615634752210745
$ws.Range("D9").Value = 8.927691427892182
$ws.Range("E9").Value = 12.86258470817622
$ws.Range("F9").Value = 32.87095191818371
$ws.Range("G9").Value = 3.622968175018213
$ws.Range("I9").Value = 17.83664119037646
$ws.Range("J9").Value = 9.57456114208212
$ws.Range("M9").Value = 21.16958263590947
$ws.Range("N9").Value = 17.35803907625224
$ws.Range("O9").Value = 24.86254992116526
$ws.Range("B10").Value = 7.954728715767244
$ws.Range("D10").Value = 8.891718100750751
$ws.Range("E10").Value = 12.77503896262109
$ws.Range("F10").Value = 33.25721904214974
$ws.Range("G10").Value = 3.61831364632175
$ws.Range("I10").Value = 17.78912326377343
$ws.Range("J10").Value = 9.517410455537222
$ws.Range("M10").Value = 21.97551605710269
$ws.Range("N10").Value = 17.46707904711397
$ws.Range("O10").Value = 25.15598466247187
$ws.Range("B11").Value = 8.104507595729507
$ws.Range("D11").Value = 8.876318341781074
$ws.Range("E11").Value = 12.73723733823538
$ws.Range("F11").Value = 33.44149579044595
$ws.Range("G11").Value = 3.616296670960564
$ws.Range("I11").Value = 17.77143094423827
$ws.Range("J11").Value = 9.492619331594808
$ws.Range("M11").Value = 22.33382761597202
$ws.Range("N11").Value = 17.51845297008461
$ws.Range("O11").Value = 25.2956434414589
$ws.Range("B12").Value = 8.160533989627812
$ws.Range("D12").Value = 8.870625216975286
$ws.Range("E12").Value = 12.72321313285768
$ws.Range("F12").Value = 33.51245750285076
$ws.Range("G12").Value = 3.615547241436046
$ws.Range("I12").Value = 17.76529421777104
$ws.Range("J12").Value = 9.483404450506432
$ws.Range("M12").Value = 22.46819910065167
$ws.Range("N12").Value = 17.53815067718478
$ws.Range("O12").Value = 25.34937863454483
$ws.Range("B13").Value = 8.148499175383952
$ws.Range("D13").Value = 8.871845181134656
$ws.Range("E13").Value = 12.72622058644522
$ws.Range("F13").Value = 33.49712298263992
$ws.Range("G13").Value = 3.615708007324851
$ws.Range("I13").Value = 17.76659085540858
$ws.Range("J13").Value = 9.485381355949727
$ws.Range("M13").Value = 22.43931999076738
$ws.Range("N13").Value = 17.53389777347758
$ws.Range("O13").Value = 25.3377686492203
$ws.Range("B14").Value = 8.109131039385737
$ws.Range("D14").Value = 8.875847192282912
$ws.Range("E14").Value = 12.73607774021069
$ws.Range("F14").Value = 33.44731050361459
$ws.Range("G14").Value = 3.616234727704544
$ws.Range("I14").Value = 17.77091479717366
$ws.Range("J14").Value = 9.491857754340863
$ws.Range("M14").Value = 22.34490940379875
$ws.Range("N14").Value = 17.52006869670069
$ws.Range("O14").Value = 25.30004746992035
$ws.Range("B15").Value = 8.084925510095299
$ws.Range("D15").Value = 8.878316558775106
$ws.Range("E15").Value = 12.74215334189713
$ws.Range("F15").Value = 33.41695104768853
$ws.Range("G15").Value = 3.616559226225911
$ws.Range("I15").Value = 17.77363661113879
$ws.Range("J15").Value = 9.495847245149255
$ws.Range("M15").Value = 22.28690581073572
$ws.Range("N15").Value = 17.51162936533273
$ws.Range("O15").Value = 25.27705161278785
$ws.Range("B16").Value = 7.944846030415742
$ws.Range("D16").Value = 8.892743898150494
$ws.Range("E16").Value = 12.77755005903375
$ws.Range("F16").Value = 33.24534410230797
$ws.Range("G16").Value = 3.618447473623847
$ws.Range("I16").Value = 17.79035835120864
$ws.Range("J16").Value = 9.519054851496209
$ws.Range("M16").Value = 21.95192199490171
$ws.Range("N16").Value = 17.4637563310027
$ws.Range("O16").Value = 25.14697862968685
$ws.Range("B17").Value = 7.857729642908057
$ws.Range("D17").Value = 8.9018414738565
$ws.Range("E17").Value = 12.79978266293179
$ws.Range("F17").Value = 33.14222618726501
$ws.Range("G17").Value = 3.619631507430579
$ws.Range("I17").Value = 17.80162073355755
$ws.Range("J17").Value = 9.533600689491315
$ws.Range("M17").Value = 21.74420449532179
$ws.Range("N17").Value = 17.43483353533229
$ws.Range("O17").Value = 25.06873763727773
$ws.Range("B18").Value = 7.80720378851348
$ws.Range("D18").Value = 8.907164976781385
$ws.Range("E18").Value = 12.81276074720886
$ws.Range("F18").Value = 33.08372460795376
$ws.Range("G18").Value = 3.62032198665322
$ws.Range("I18").Value = 17.8084679737602
$ws.Range("J18").Value = 9.542080714758896
$ws.Range("M18").Value = 21.62395503151471
$ws.Range("N18").Value = 17.41836529728024
$ws.Range("O18").Value = 25.02431901779956
$ws.Range("B19").Value = 7.790026207901633
$ws.Range("D19").Value = 8.908983030462265
$ws.Range("E19").Value = 12.81718763432826
$ws.Range("F19").Value = 33.06405745923142
$ws.Range("G19").Value = 3.620557397317653
$ws.Range("I19").Value = 17.81084981499088
$ws.Range("J19").Value = 9.54497144211317
$ws.Range("M19").Value = 21.58311112648971
$ws.Range("N19").Value = 17.41281851990801
$ws.Range("O19").Value = 25.00938099670181
$ws.Range("B20").Value = 7.867047038460559
$ws.Range("D20").Value = 8.900863624207851
$ws.Range("E20").Value = 12.7973962542917
$ws.Range("F20").Value = 33.15311988689148
$ws.Range("G20").Value = 3.61950448711842
$ws.Range("I20").Value = 17.8003836092916
$ws.Range("J20").Value = 9.532040502801111
$ws.Range("M20").Value = 21.76639756053003
$ws.Range("N20").Value = 17.43789517224297
$ws.Range("O20").Value = 25.0770064036021
$ws.Range("B21").Value = 8.120713543693119
$ws.Range("D21").Value = 8.874667949943827
$ws.Range("E21").Value = 12.73317457643282
$ws.Range("F21").Value = 33.46191002789835
$ws.Range("G21").Value = 3.616079628207121
$ws.Range("I21").Value = 17.76962948277678
$ws.Range("J21").Value = 9.489950789155355
$ws.Range("M21").Value = 22.37267662080097
$ws.Range("N21").Value = 17.52412411301059
$ws.Range("O21").Value = 25.31110435531908
$ws.Range("B22").Value = 8.2824490644402
$ws.Range("D22").Value = 8.858354306455128
$ws.Range("E22").Value = 12.69289473992688
$ws.Range("F22").Value = 33.67057541715511
$ws.Range("G22").Value = 3.613924921134488
$ws.Range("I22").Value = 17.75281059623247
$ws.Range("J22").Value = 9.463450720129847
$ws.Range("M22").Value = 22.76121163545895
$ws.Range("N22").Value = 17.58189376384235
$ws.Range("O22").Value = 25.46903313325213
$ws.Range("B23").Value = 8.196512976363133
$ws.Range("D23").Value = 8.8669874823106
$ws.Range("E23").Value = 12.71423812256202
$ws.Range("F23").Value = 33.55859721086604
$ws.Range("G23").Value = 3.615067302810445
$ws.Range("I23").Value = 17.76148743171519
$ws.Range("J23").Value = 9.477502265641654
$ws.Range("M23").Value = 22.55458470092979
$ws.Range("N23").Value = 17.55093538064273
$ws.Range("O23").Value = 25.38430538481069
$ws.Range("B24").Value = 7.862836015368699
$ws.Range("D24").Value = 8.90130541975995
$ws.Range("E24").Value = 12.79847453826342
$ws.Range("F24").Value = 33.14819240468627
$ws.Range("G24").Value = 3.619561882586382
$ws.Range("I24").Value = 17.80094175328543
$ws.Range("J24").Value = 9.532745497378672
$ws.Range("M24").Value = 21.7563666516248
$ws.Range("N24").Value = 17.43651050740555
$ws.Range("O24").Value = 25.07326633721634
$ws.Range("B25").Value = 7.487135860670749
$ws.Range("D25").Value = 8.941796338785757
$ws.Range("E25").Value = 12.89662200886409
$ws.Range("F25").Value = 32.73661863382011
$ws.Range("G25").Value = 3.624771815638837
$ws.Range("I25").Value = 17.85762662433649
$ws.Range("J25").Value = 9.596681714212338
$ws.Range("M25").Value = 20.86625575037261
$ws.Range("N25").Value = 17.31952653314243
$ws.Range("O25").Value = 24.76020610028862
